$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 2956
$ws.Range("F7").Value = 2313
$ws.Range("F8").Value = 1646
$ws.Range("F9").Value = 53
$ws.Range("F11").Value = 108
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 2651
$ws.Range("F15").Value = 1503
$ws.Range("F16").Value = 7007
$ws.Range("F18").Value = 7169
$ws.Range("F20").Value = 5367
$ws.Range("F21").Value = 3095
$ws.Range("F22").Value = 3466
$ws.Range("F24").Value = 167
$ws.Range("F25").Value = 1859
$ws.Range("F26").Value = 77
$ws.Range("F27").Value = 296
$ws.Range("F28").Value = 874
$ws.Range("F30").Value = 173
$ws.Range("F31").Value = 37
$ws.Range("F32").Value = 2394
$ws.Range("F33").Value = 1147
$ws.Range("F34").Value = 2626
$ws.Range("F35").Value = 16
$ws.Range("F36").Value = 18
$ws.Range("F38").Value = 376
$ws.Range("F39").Value = 1044
$ws.Range("F41").Value = 469
$ws.Range("F42").Value = 519

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 211
$ws.Range("F9").Value = 32
$ws.Range("F14").Value = 92

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 2956
$ws.Range("F8").Value = 2313
$ws.Range("F9").Value = 1646
$ws.Range("F10").Value = 53
$ws.Range("F12").Value = 108
$ws.Range("F14").Value = 2651
$ws.Range("F15").Value = 1503
$ws.Range("F16").Value = 211
$ws.Range("F17").Value = 32
$ws.Range("F19").Value = 7007
$ws.Range("F21").Value = 7169
$ws.Range("F23").Value = 5367
$ws.Range("F24").Value = 3095
$ws.Range("F25").Value = 3466
$ws.Range("F29").Value = 1859
$ws.Range("F32").Value = 296
$ws.Range("F33").Value = 874
$ws.Range("F35").Value = 173
$ws.Range("F36").Value = 37
$ws.Range("F37").Value = 2394
$ws.Range("F38").Value = 1147
$ws.Range("F40").Value = 2626
$ws.Range("F41").Value = 16
$ws.Range("F42").Value = 18
$ws.Range("F45").Value = 376
$ws.Range("F46").Value = 1044
$ws.Range("F48").Value = 469
$ws.Range("F49").Value = 519
